# Weekly data refresh: a new price observation is inserted as row 202
# (shifting the existing rows 202-229 down to 203-230) on the single
# worksheet of the "Berenjena" workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 202; everything currently at/after row 202 moves
# down one row (old 202 -> 203, ..., old 229 -> 230).
$ws.Rows(202).Insert()

# Populate the newly inserted row 202 with the new weekly record.
$ws.Cells.Item(202, 1).Value = 6
$ws.Cells.Item(202, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(202, 3).Value = "Metropolitana"
$ws.Cells.Item(202, 4).Value = 44776
$ws.Cells.Item(202, 5).Value = 13
$ws.Cells.Item(202, 6).Value = 100112001
$ws.Cells.Item(202, 7).Value = "Berenjena"
$ws.Cells.Item(202, 8).Value = "Sin especificar"
$ws.Cells.Item(202, 9).Value = "Primera"
$ws.Cells.Item(202, 10).Value = 650
$ws.Cells.Item(202, 11).Value = 10000
$ws.Cells.Item(202, 12).Value = 11000
$ws.Cells.Item(202, 13).Value = 10646
$ws.Cells.Item(202, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(202, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(202, 16).Value = 213
$ws.Cells.Item(202, 17).Value = 50
$ws.Cells.Item(202, 18).Value = "Hortaliza"

# Apply the same date number format (style index 2 in the original file,
# "YYYY-MM-DD HH:MM:SS") used by the other cells in the Fecha column so the
# new row's date cell is formatted consistently.
$ws.Cells.Item(202, 4).NumberFormat = $ws.Cells.Item(203, 4).NumberFormat
